$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.560.11"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.577.45"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3691"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.62"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3341"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.146"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07473"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.005"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.964"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "1.576.97"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001117"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06763"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.428"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "22.545.51"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.402"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.606"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.018"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("D31").Value = "1.753.29"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.074"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.185"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.006"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.700"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08324"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02457"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.441"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.302"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06389"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6354"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6206"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.94%  "
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.220"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07270"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.85%  "
